$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2013888888888889
$ws.Range("C2").Value = 0.5694444444444444
$ws.Range("J2").Value = 0.003472222222222222
$ws.Range("P2").Value = 0.1423611111111111
$ws.Range("S2").Value = 0.08333333333333333
$ws.Range("C3").Value = 0.03012048192771084
$ws.Range("J3").Value = 0.01807228915662651
$ws.Range("P3").Value = 0.7228915662650602
$ws.Range("S3").Value = 0.2289156626506024
$ws.Range("J4").Value = 0.1282051282051282
$ws.Range("P4").Value = 0.5128205128205128
$ws.Range("S4").Value = 0.358974358974359
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.05181347150259067
$ws.Range("D6").Value = 0.005181347150259068
$ws.Range("F6").Value = 0.05699481865284974
$ws.Range("J6").Value = 0.1761658031088083
$ws.Range("O6").Value = 0.02590673575129534
$ws.Range("Q6").Value = 0.1398963730569948
$ws.Range("R6").Value = 0.09326424870466321
$ws.Range("S6").Value = 0.4507772020725389
$ws.Range("B7").Value = 0.1575342465753425
$ws.Range("D7").Value = 0.0136986301369863
$ws.Range("F7").Value = 0.0684931506849315
$ws.Range("J7").Value = 0.08904109589041095
$ws.Range("O7").Value = 0.0136986301369863
$ws.Range("Q7").Value = 0.1849315068493151
$ws.Range("R7").Value = 0.0547945205479452
$ws.Range("S7").Value = 0.4178082191780822
$ws.Range("B8").Value = 0.1128608923884514
$ws.Range("D8").Value = 0.01312335958005249
$ws.Range("F8").Value = 0.05774278215223097
$ws.Range("J8").Value = 0.09711286089238845
$ws.Range("O8").Value = 0.005249343832020997
$ws.Range("Q8").Value = 0.1627296587926509
$ws.Range("R8").Value = 0.08923884514435695
$ws.Range("S8").Value = 0.4619422572178478
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.01904761904761905
$ws.Range("F9").Value = 0.09523809523809523
$ws.Range("J9").Value = 0.1285714285714286
$ws.Range("O9").Value = 0.009523809523809525
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.08095238095238096
$ws.Range("S9").Value = 0.3714285714285714
$ws.Range("B10").Value = 0.1154529307282416
$ws.Range("D10").Value = 0.02486678507992895
$ws.Range("E10").Value = 0.003552397868561279
$ws.Range("F10").Value = 0.07104795737122557
$ws.Range("J10").Value = 0.1136767317939609
$ws.Range("O10").Value = 0.009769094138543518
$ws.Range("Q10").Value = 0.2015985790408526
$ws.Range("R10").Value = 0.08880994671403197
$ws.Range("S10").Value = 0.3712255772646537
$ws.Range("G11").Value = 0.1769911504424779
$ws.Range("J11").Value = 0.09734513274336283
$ws.Range("K11").Value = 0.2300884955752212
$ws.Range("L11").Value = 0.4867256637168141
$ws.Range("S11").Value = 0.008849557522123894
$ws.Range("G12").Value = 0.7280701754385965
$ws.Range("J12").Value = 0.1754385964912281
$ws.Range("K12").Value = 0.03508771929824561
$ws.Range("L12").Value = 0.02631578947368421
$ws.Range("S12").Value = 0.03508771929824561
$ws.Range("G13").Value = 0.6222222222222222
$ws.Range("J13").Value = 0.3111111111111111
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("H15").Value = 0.1693121693121693
$ws.Range("I15").Value = 0.07936507936507936
$ws.Range("J15").Value = 0.4338624338624338
$ws.Range("K15").Value = 0.07936507936507936
$ws.Range("M15").Value = 0.01058201058201058
$ws.Range("O15").Value = 0.05291005291005291
$ws.Range("S15").Value = 0.1746031746031746
$ws.Range("F16").Value = 0.005780346820809248
$ws.Range("H16").Value = 0.138728323699422
$ws.Range("I16").Value = 0.1213872832369942
$ws.Range("J16").Value = 0.4508670520231214
$ws.Range("K16").Value = 0.04624277456647399
$ws.Range("M16").Value = 0.02890173410404624
$ws.Range("O16").Value = 0.05202312138728324
$ws.Range("S16").Value = 0.1560693641618497
$ws.Range("F17").Value = 0.02590673575129534
$ws.Range("H17").Value = 0.1658031088082902
$ws.Range("I17").Value = 0.09326424870466321
$ws.Range("J17").Value = 0.422279792746114
$ws.Range("K17").Value = 0.05958549222797927
$ws.Range("M17").Value = 0.0155440414507772
$ws.Range("O17").Value = 0.05440414507772021
$ws.Range("S17").Value = 0.1632124352331606
$ws.Range("F18").Value = 0.02259887005649718
$ws.Range("H18").Value = 0.1807909604519774
$ws.Range("I18").Value = 0.096045197740113
$ws.Range("J18").Value = 0.4124293785310734
$ws.Range("K18").Value = 0.0903954802259887
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.03389830508474576
$ws.Range("S18").Value = 0.1468926553672316
$ws.Range("F19").Value = 0.005857740585774059
$ws.Range("H19").Value = 0.1907949790794979
$ws.Range("I19").Value = 0.09958158995815899
$ws.Range("J19").Value = 0.3707112970711297
$ws.Range("K19").Value = 0.08870292887029289
$ws.Range("M19").Value = 0.02594142259414226
$ws.Range("N19").Value = 0.0008368200836820083
$ws.Range("O19").Value = 0.08702928870292886
$ws.Range("S19").Value = 0.1305439330543933
